$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price cells to remain text so values that look like
# numbers (e.g. "551.97", "1.00") are not auto-converted by Excel.
$textCells = @(
    "D5", "D6", "D7", "D8", "D10", "D12", "D14", "D17", "D20", "D21", "D22", "D23", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.168.20"
$ws.Range("E2").Value = "  +0.49%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.319.59"
$ws.Range("E3").Value = "  -0.20%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "551.97"
$ws.Range("E5").Value = "  +0.18%  "

# Row 6 - Solana
$ws.Range("D6").Value = "172.87"
$ws.Range("E6").Value = "  +0.57%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.617"
$ws.Range("E7").Value = "  +0.86%  "

# Row 8 - USDC
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.309.88"

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +5.82%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +1.26%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "53.27"
$ws.Range("E12").Value = "  +0.81%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +3.52%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "9.03"
$ws.Range("E14").Value = "  +0.71%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.850.74"
$ws.Range("E15").Value = "  -0.38%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +3.02%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "18.06"
$ws.Range("E17").Value = "  -0.77%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.308.88"
$ws.Range("E18").Value = "  -0.63%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "64.065.97"
$ws.Range("E19").Value = "  +0.30%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "11.64"
$ws.Range("E20").Value = "  -1.05%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "0.981"
$ws.Range("E21").Value = "  +1.35%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "455.36"
$ws.Range("E22").Value = "  +6.99%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "5.08"
$ws.Range("E23").Value = "  +8.90%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  -0.53%  "

# Row 25 - InternetComputer(DFINITY) (was Litecoin)
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "13.84"
$ws.Range("E25").Value = "  +3.94%  "

# Row 26 - Litecoin (was InternetComputer(DFINITY))
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "86.60"
$ws.Range("E26").Value = "  +3.09%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  +1.37%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "10.60"
$ws.Range("E28").Value = "  +0.10%  "

# Row 29 - EthereumClassic (was Filecoin)
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "30.88"
$ws.Range("E29").Value = "  +4.53%  "

# Row 30 - Filecoin (was EthereumClassic)
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "8.54"
$ws.Range("E30").Value = "  +0.11%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "6.52"
$ws.Range("E31").Value = "  -2.10%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "11.38"

# Row 33 - OKB
$ws.Range("D33").Value = "61.50"
$ws.Range("E33").Value = "  +5.92%  "

# Row 34 - Bittensor
$ws.Range("D34").Value = "564.51"
$ws.Range("E34").Value = "  -4.95%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  -0.05%  "

# Row 36 - Dai
$ws.Range("E36").Value = "  +0.10%  "

# Row 37 - Stacks (was Kaspa)
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "3.53"
$ws.Range("E37").Value = "  +2.19%  "

# Row 38 - Kaspa (was Stacks)
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  -1.96%  "

# Row 39 - InjectiveProtocol
$ws.Range("D39").Value = "35.12"
$ws.Range("E39").Value = "  -0.07%  "

# Row 40 - TheGraph
$ws.Range("E40").Value = "  +0.36%  "

# Row 41 - PEPE
$ws.Range("E41").Value = "  -2.32%  "

# Row 42 - Maker
$ws.Range("D42").Value = "3.042.53"
$ws.Range("E42").Value = "  -1.33%  "

# Row 43 - VeChain
$ws.Range("D43").Value = "0.0415"
$ws.Range("E43").Value = "  +2.63%  "

# Row 44 - ThetaToken
$ws.Range("D44").Value = "2.74"
$ws.Range("E44").Value = "  -1.15%  "

# Row 45 - ApeXProtocol
$ws.Range("D45").Value = "3.21"
$ws.Range("E45").Value = "  +0.90%  "

# Row 46 - Fetch.AI
$ws.Range("D46").Value = "2.44"
$ws.Range("E46").Value = "  +0.68%  "

# Row 47 - Stellar
$ws.Range("D47").Value = "0.133"
$ws.Range("E47").Value = "  +2.88%  "

# Row 48 - FirstDigitalUSD
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  -0.08%  "

# Row 49 - Monero
$ws.Range("D49").Value = "140.70"
$ws.Range("E49").Value = "  +6.45%  "

# Row 50 - WEMIXToken
$ws.Range("D50").Value = "2.49"
$ws.Range("E50").Value = "  -3.58%  "

# Row 51 - THORChain
$ws.Range("D51").Value = "8.11"
$ws.Range("E51").Value = "  -0.29%  "

